$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-08 13:33:22"
$wsZhCn.Range("G3").Value = "2016-01-08 13:34:12"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-08 13:33:35"
$wsDeDe.Range("G3").Value = "2016-01-08 13:34:33"
